$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Apply formatting first (style 2 = header row, style 3 = data rows) ---
# Header row 1: B1:G1 should carry the same style as A1 (style 2).
$ws.Range("A1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Data row 2: B2:G2 should carry the same style as A2 (style 3).
$ws.Range("A2").Copy()
$ws.Range("B2:G2").PasteSpecial(-4122)

# Data row 3: B3,D3:G3 should carry style 3 -- C3 is intentionally skipped/untouched.
$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("D3:G3").PasteSpecial(-4122)

# Data row 4 (new row): D4:G4 carry style 3, B4 stays default (no explicit style), A4/C4 stay empty.
$ws.Range("A3").Copy()
$ws.Range("D4:G4").PasteSpecial(-4122)

# --- 2. Write the final cell values/content, in the order the strings were first
#        authored so the shared-string table comes out in the same sequence. ---

# Header row.
$ws.Range("A1").Value = "ID (optional)"
$ws.Range("B1").Value = "en_comments"
$ws.Range("C1").Value = "de_comments"
$ws.Range("D1").Value = "en_list"
$ws.Range("E1").Value = "de_list"
$ws.Range("F1").Value = "en_1"
$ws.Range("G1").Value = "de_1"

# Pre-existing rows keep their list/node values, shifted into the new columns.
$ws.Range("B2").ClearContents()
$ws.Range("D2").Value = "List 4"
$ws.Range("E2").Value = "List de"

$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C3").ClearFormats()
$ws.Range("D3").Value = "List 4"
$ws.Range("E3").Value = "List de"
$ws.Range("F3").Value = "Node 1"

# New row 4 (Node 2 / Knoten 2) filled in first...
$ws.Range("D4").Value = "List 4"
$ws.Range("E4").Value = "List de"
$ws.Range("F4").Value = "Node 2"
$ws.Range("G4").Value = "Knoten 2"

# ...then the two "missing translation" comment cells are filled in afterwards.
$ws.Range("C2").Value = "List Missing English Comment"
$ws.Range("B4").Value = "Node Missing German Comment"

Write-Host "done"
